# Updated symbol list on Fri Jan 27 07:21:51 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for the crypto symbol table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store their numeric-looking figures as literal text
# (trailing zeros, "--", percent signs, etc. must survive verbatim), so force
# the target range to Text format before writing the new values. This mirrors
# how the sheet was originally populated and stops Excel from silently
# re-interpreting e.g. "304.55" or "-0.91%" as a Number/Percentage.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "304.55"
$ws.Range("E2").Value = "-0.91%"
$ws.Range("D3").Value = "35.74"
$ws.Range("E3").Value = "-0.40%"
$ws.Range("D4").Value = "5.039"
$ws.Range("E4").Value = "-1.38%"
$ws.Range("D5").Value = "0.07952"
$ws.Range("E5").Value = "-1.67%"
$ws.Range("D6").Value = "1.932"
$ws.Range("E6").Value = "-0.97%"
$ws.Range("D7").Value = "7.778"
$ws.Range("E7").Value = "0.24%"
$ws.Range("D8").Value = "0.9192"
$ws.Range("E8").Value = "-1.11%"
$ws.Range("D9").Value = "0.1308"
$ws.Range("E9").Value = "-5.40%"
$ws.Range("D10").Value = "0.1914"
$ws.Range("E10").Value = "-0.59%"
$ws.Range("D11").Value = "0.09086"
$ws.Range("E11").Value = "-1.53%"
$ws.Range("D12").Value = "0.03440"
$ws.Range("E12").Value = "-0.38%"
$ws.Range("D13").Value = "0.09844"
$ws.Range("E13").Value = "0.10%"
$ws.Range("D14").Value = "0.001404"
$ws.Range("E14").Value = "-1.28%"
$ws.Range("D15").Value = "0.006123"
$ws.Range("E15").Value = "6.55%"
$ws.Range("D16").Value = "3.727"
$ws.Range("E16").Value = "2.93%"
$ws.Range("D17").Value = "4.125"
$ws.Range("E17").Value = "-1.94%"
$ws.Range("D18").Value = "3.370"
$ws.Range("E18").Value = "13.49%"
$ws.Range("D19").Value = "0.3444"
$ws.Range("E19").Value = "0.16%"
$ws.Range("D20").Value = "0.1311"
$ws.Range("E20").Value = "-2.22%"
$ws.Range("D21").Value = "5.166"
$ws.Range("E21").Value = "5.50%"
$ws.Range("D22").Value = "0.2351"
$ws.Range("E22").Value = "-3.80%"
$ws.Range("D23").Value = "0.04413"
$ws.Range("E23").Value = "-0.93%"
$ws.Range("D24").Value = "0.001232"
$ws.Range("E24").Value = "0.93%"
$ws.Range("D25").Value = "0.004625"
$ws.Range("E25").Value = "-4.34%"
$ws.Range("E26").Value = "0.56%"
$ws.Range("D27").Value = "0.0004442"
$ws.Range("E27").Value = "0.01%"
$ws.Range("D39").Value = "0.01936"
$ws.Range("E39").Value = "-4.34%"
$ws.Range("D40").Value = "0.05359"
$ws.Range("E40").Value = "8.58%"
$ws.Range("D41").Value = "0.007638"
$ws.Range("E41").Value = "-0.92%"
$ws.Range("E42").Value = "0.46%"
$ws.Range("D43").Value = "0.1352"
$ws.Range("E43").Value = "-1.91%"
$ws.Range("D44").Value = "0.002140"
$ws.Range("E44").Value = "1.66%"
$ws.Range("D45").Value = "0.009579"
$ws.Range("E45").Value = "-17.51%"
$ws.Range("D46").Value = "0.00006162"
$ws.Range("E46").Value = "-4.45%"
$ws.Range("E47").Value = "-0.12%"
$ws.Range("D48").Value = "63.57"
$ws.Range("E48").Value = "-1.69%"
$ws.Range("E49").Value = "39.11%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "-0.12%"
$ws.Range("E51").Value = "-0.12%"
